# Update the "Authors" column (E2:E16) with corrected author-list strings.
# (commit: "cleaned more data, fixed a problem where some SLRs contained incorrect data.")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[Barnaby%D.P.%coreGivesNoEmail%0,              Becker%L.B.%coreGivesNoEmail%0,              Chelico%J.D.%coreGivesNoEmail%0,              Cohen%S.L.%coreGivesNoEmail%0,              Cookingham%J.%coreGivesNoEmail%0,              Coppa%K.%coreGivesNoEmail%0,              Crawford%J.M.%coreGivesNoEmail%0,              Davidson%K.W.%coreGivesNoEmail%0,              Diefenbach%M.A.%coreGivesNoEmail%0,              Dominello%A.J.%coreGivesNoEmail%0,              Duer-Hefele%J.%coreGivesNoEmail%0,              Falzon%L.%coreGivesNoEmail%0,              Gitlin%J.%coreGivesNoEmail%0,              Hajizadeh%N.%coreGivesNoEmail%0,              Harvin%T.G.%coreGivesNoEmail%0,              Hirsch%J.S.%coreGivesNoEmail%0,              Hirschwerk%D.A.%coreGivesNoEmail%0,              Kim%E.J.%coreGivesNoEmail%0,              Kozel%Z.M.%coreGivesNoEmail%0,              Marrast%L.M.%coreGivesNoEmail%0,              McGinn%T.%coreGivesNoEmail%0,              Mogavero%J.N.%coreGivesNoEmail%0,              Narasimhan%M.%coreGivesNoEmail%0,              Osorio%G.A.%coreGivesNoEmail%0,              Qiu%M.%coreGivesNoEmail%0,              Richardson%S.%coreGivesNoEmail%0,              Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Range("E3").Value2 = '[Kai%Liu%hmliukai@126.com%1,                 Ying%Chen%NULL%1,                 Ruzheng%Lin%NULL%1,                 Kunyuan%Han%NULL%1]'
$ws.Range("E4").Value2 = '[Pingzheng%Mo%NULL%4,                 Yuanyuan%Xing%NULL%1,                 Yu%Xiao%NULL%0,                 Liping%Deng%NULL%1,                 Qiu%Zhao%NULL%1,                 Hongling%Wang%NULL%1,                 Yong%Xiong%NULL%3,                 Zhenshun%Cheng%NULL%4,                 Shicheng%Gao%NULL%1,                 Ke%Liang%NULL%1,                 Mingqi%Luo%NULL%1,                 Tielong%Chen%NULL%1,                 Shihui%Song%NULL%1,                 Zhiyong%Ma%NULL%1,                 Xiaoping%Chen%NULL%1,                 Ruiying%Zheng%NULL%1,                 Qian%Cao%NULL%1,                 Fan%Wang%fanndywang@foxmail.com%0,                 Yongxi%Zhang%znact1936@126.com%0]'
$ws.Range("E5").Value2 = '[Carly%Eastin%NULL%0,                 Travis%Eastin%NULL%1]'
$ws.Range("E6").Value2 = '[Tao%Guo%NULL%0,                 Yongzhen%Fan%NULL%2,                 Ming%Chen%NULL%2,                 Xiaoyan%Wu%NULL%2,                 Lin%Zhang%NULL%2,                 Tao%He%NULL%2,                 Hairong%Wang%NULL%2,                 Jing%Wan%NULL%2,                 Xinghuan%Wang%NULL%0,                 Zhibing%Lu%NULL%0]'
$ws.Range("E7").Value2 = '[Chaolin%Huang%NULL%0,                 Yeming%Wang%NULL%0,                 Xingwang%Li%NULL%4,                 Lili%Ren%NULL%3,                 Jianping%Zhao%NULL%10,                 Yi%Hu%NULL%5,                 Li%Zhang%NULL%5,                 Guohui%Fan%NULL%0,                 Jiuyang%Xu%NULL%0,                 Xiaoying%Gu%NULL%0,                 Zhenshun%Cheng%NULL%0,                 Ting%Yu%NULL%0,                 Jiaan%Xia%NULL%3,                 Yuan%Wei%NULL%0,                 Wenjuan%Wu%NULL%3,                 Xuelei%Xie%NULL%3,                 Wen%Yin%NULL%3,                 Hui%Li%NULL%0,                 Min%Liu%NULL%4,                 Yan%Xiao%NULL%4,                 Hong%Gao%NULL%5,                 Li%Guo%NULL%4,                 Jungang%Xie%NULL%5,                 Guangfa%Wang%NULL%3,                 Rongmeng%Jiang%NULL%3,                 Zhancheng%Gao%NULL%3,                 Qi%Jin%NULL%4,                 Jianwei%Wang%wangjw28@163.com%3,                 Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("E8").Value2 = '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,              Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%0,              Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%0,              Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%0,              Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%0]'
$ws.Range("E9").Value2 = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,              Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%0,              Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%0,              Subert-Salas%Lizandra%coreGivesNoEmail%0]'
$ws.Range("E10").Value2 = '[Antonelli%Massimo%coreGivesNoEmail%7,              Cabrini%Luca%coreGivesNoEmail%2,              Castelli%Antonio%coreGivesNoEmail%2,              Cecconi%Maurizio%coreGivesNoEmail%2,              Cereda%Danilo%coreGivesNoEmail%2,              Coluccello%Antonio%coreGivesNoEmail%2,              Foti%Giuseppe%coreGivesNoEmail%2,              Fumagalli%Roberto%coreGivesNoEmail%2,              Grasselli%Giacomo%coreGivesNoEmail%2,              Iotti%Giorgio%coreGivesNoEmail%2,              Latronico%Nicola%coreGivesNoEmail%2,              Lorini%Luca%coreGivesNoEmail%2,              Merler%Stefano%coreGivesNoEmail%2,              Natalini%Giuseppe%coreGivesNoEmail%2,              Pesenti%Antonio%coreGivesNoEmail%2,              Piatti%Alessandra%coreGivesNoEmail%2,              Ranieri%Marco Vito%coreGivesNoEmail%2,              Scandroglio%Anna Mara%coreGivesNoEmail%2,              Storti%Enrico%coreGivesNoEmail%2,              Zanella%Alberto%coreGivesNoEmail%2,              Zangrillo%Alberto%coreGivesNoEmail%2]'
$ws.Range("E11").Value2 = '[Pavan K.%Bhatraju%NULL%0,                 Bijan J.%Ghassemieh%NULL%4,                 Michelle%Nichols%NULL%4,                 Richard%Kim%NULL%4,                 Keith R.%Jerome%NULL%4,                 Arun K.%Nalla%NULL%4,                 Alexander L.%Greninger%NULL%4,                 Sudhakar%Pipavath%NULL%4,                 Mark M.%Wurfel%NULL%4,                 Laura%Evans%NULL%4,                 Patricia A.%Kritek%NULL%4,                 T. Eoin%West%NULL%4,                 Andrew%Luks%NULL%4,                 Anthony%Gerbino%NULL%4,                 Chris R.%Dale%NULL%4,                 Jason D.%Goldman%NULL%4,                 Shane%O’Mahony%NULL%4,                 Carmen%Mikacenic%NULL%4]'
$ws.Range("E12").Value2 = '[Matt%Arentz%xref no email%7,          Eric%Yim%xref no email%7,          Lindy%Klaff%xref no email%7,          Sharukh%Lokhandwala%xref no email%7,          Francis X.%Riedo%xref no email%7,          Maria%Chong%xref no email%7,          Melissa%Lee%xref no email%7]'
$ws.Range("E13").Value2 = '[Chaomin%Wu%xref no email%0,          Xiaoyan%Chen%xref no email%0,          Yanping%Cai%xref no email%0,          Jia\u2019an%Xia%xref no email%0,          Xing%Zhou%xref no email%0,          Sha%Xu%xref no email%0,          Hanping%Huang%xref no email%0,          Li%Zhang%xref no email%0,          Xia%Zhou%xref no email%0,          Chunling%Du%xref no email%0,          Yuye%Zhang%xref no email%0,          Juan%Song%xref no email%0,          Sijiao%Wang%xref no email%0,          Yencheng%Chao%xref no email%0,          Zeyong%Yang%xref no email%0,          Jie%Xu%xref no email%0,          Xin%Zhou%xref no email%0,          Dechang%Chen%xref no email%0,          Weining%Xiong%xref no email%0,          Lei%Xu%xref no email%0,          Feng%Zhou%xref no email%0,          Jinjun%Jiang%xref no email%0,          Chunxue%Bai%xref no email%0,          Junhua%Zheng%xref no email%0,          Yuanlin%Song%xref no email%0]'
$ws.Range("E14").Value2 = '[Nanshan%Chen%NULL%0,                 Min%Zhou%NULL%4,                 Xuan%Dong%NULL%1,                 Jieming%Qu%NULL%4,                 Fengyun%Gong%NULL%1,                 Yang%Han%NULL%1,                 Yang%Qiu%NULL%1,                 Jingli%Wang%NULL%1,                 Ying%Liu%NULL%0,                 Yuan%Wei%NULL%0,                 Jia''an%Xia%NULL%6,                 Ting%Yu%NULL%0,                 Xinxin%Zhang%NULL%1,                 Li%Zhang%NULL%0]'
$ws.Range("E15").Value2 = '[Fei%Zhou%NULL%0,                 Ting%Yu%NULL%0,                 Ronghui%Du%NULL%0,                 Guohui%Fan%NULL%0,                 Ying%Liu%NULL%0,                 Zhibo%Liu%NULL%0,                 Jie%Xiang%NULL%0,                 Yeming%Wang%NULL%0,                 Bin%Song%NULL%0,                 Xiaoying%Gu%NULL%0,                 Lulu%Guan%NULL%0,                 Yuan%Wei%NULL%0,                 Hui%Li%NULL%0,                 Xudong%Wu%NULL%0,                 Jiuyang%Xu%NULL%0,                 Shengjin%Tu%NULL%0,                 Yi%Zhang%NULL%0,                 Hua%Chen%NULL%0,                 Bin%Cao%NULL%0]'
$ws.Range("E16").Value2 = '[Ke%Wang%NULL%0,                 Wei%Zhao%NULL%1,                 Ji%Li%NULL%1,                 Weiwei%Shu%shuweiwei361@163.com%1,                 Jun%Duan%duanjun412589@163.com%1]'
